$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell updates derived from the commit diff (symbol list refresh).
# Columns D (Price) and E (Volume 1h) store numeric/percent-looking
# values as literal text, so we force a text number format before
# assigning them to avoid Excel auto-converting them to numbers.

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '259.56'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '5.91%'
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '27.86'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '-2.70%'
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '5.220'
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '-0.41%'
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.05945'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '4.29%'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '6.737'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '1.82%'
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.8727'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '2.58%'
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.9894'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '15.61%'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '4.36%'
$ws.Range('B10').Value = 'MandalaExchangeToken'
$ws.Range('C10').Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.07260'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '2.60%'
$ws.Range('B11').Value = 'BitrueCoin'
$ws.Range('C11').Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.03244'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '2.25%'
$ws.Range('B12').Value = 'BitMartToken'
$ws.Range('C12').Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.09240'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '0.41%'
$ws.Range('B13').Value = 'BitForexToken'
$ws.Range('C13').Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.001549'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '1.17%'
$ws.Range('B14').Value = 'One'
$ws.Range('C14').Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.0006038'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '0.98%'
$ws.Range('B15').Value = 'TigerCash'
$ws.Range('C15').Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.005866'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '-3.13%'
$ws.Range('B16').Value = 'LEO'
$ws.Range('C16').Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '3.496'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '0.18%'
$ws.Range('B17').Value = 'GateToken'
$ws.Range('C17').Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '3.246'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '1.84%'
$ws.Range('B18').Value = 'BTSEToken'
$ws.Range('C18').Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '2.210'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '1.65%'
$ws.Range('B19').Value = 'BitpandaEcosystemToken'
$ws.Range('C19').Value = 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.3172'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '0.16%'
$ws.Range('B20').Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range('C20').Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.03632'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '11.13%'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.1291'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '0.23%'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '3.531'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '1.37%'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.04176'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '2.54%'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '1.25%'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.001216'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '-0.38%'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.004580'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '10.69%'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.0001197'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '-0.15%'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '33.49%'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.03870'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '3.16%'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.005399'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '45.24%'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.1110'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '4.53%'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.002377'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '-4.47%'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.01093'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '16.98%'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.00005423'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '2.64%'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.00000000748'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '-0.20%'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.08532'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '13.77%'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.002135'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '-12.42%'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.00002096'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '-0.20%'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0001996'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '-0.20%'
